$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.389221
$ws.Range("H2").Value = 4.167663
$ws.Range("I2").Value = 0.2910270461264192
$ws.Range("J2").Value = 0.2910270461264192
$ws.Range("M2").Value = 8.820647333333334
$ws.Range("N2").Value = 26.461942
$ws.Range("O2").Value = 0.06415146660411865
$ws.Range("P2").Value = 0.06415146660411865
$ws.Range("Q2").Value = 12.25382850906067
$ws.Range("R2").Value = 110.284456581546
$ws.Range("S2").Value = 0.01866981183047428
$ws.Range("T2").Value = 0.01866981183047428
$ws.Range("G3").Value = 1.389221
$ws.Range("H3").Value = 4.167663
$ws.Range("I3").Value = 0.2910270461264192
$ws.Range("J3").Value = 0.2910270461264192
$ws.Range("O3").Value = 0.3979101621202897
$ws.Range("P3").Value = 0.3979101621202898
$ws.Range("Q3").Value = 76.006413364235
$ws.Range("R3").Value = 684.057720278115
$ws.Range("S3").Value = 0.1158026191055525
$ws.Range("T3").Value = 0.1158026191055525
$ws.Range("G4").Value = 1.389221
$ws.Range("H4").Value = 4.167663
$ws.Range("I4").Value = 0.2910270461264192
$ws.Range("J4").Value = 0.2910270461264192
$ws.Range("M4").Value = 21.90816333333333
$ws.Range("N4").Value = 65.72449
$ws.Range("O4").Value = 0.1593353362087987
$ws.Range("P4").Value = 0.1593353362087987
$ws.Range("Q4").Value = 30.43528057409667
$ws.Range("R4").Value = 273.91752516687
$ws.Range("S4").Value = 0.04637089224040655
$ws.Range("T4").Value = 0.04637089224040655
$ws.Range("G5").Value = 1.389221
$ws.Range("H5").Value = 4.167663
$ws.Range("I5").Value = 0.2910270461264192
$ws.Range("J5").Value = 0.2910270461264192
$ws.Range("M5").Value = 52.056859
$ws.Range("N5").Value = 156.170577
$ws.Range("O5").Value = 0.3786030350667928
$ws.Range("P5").Value = 0.3786030350667929
$ws.Range("Q5").Value = 72.318481716839
$ws.Range("R5").Value = 650.866335451551
$ws.Range("S5").Value = 0.1101837229499858
$ws.Range("T5").Value = 0.1101837229499858
$ws.Range("I6").Value = 0.461328155686921
$ws.Range("J6").Value = 0.4613281556869209
$ws.Range("M6").Value = 8.820647333333334
$ws.Range("N6").Value = 26.461942
$ws.Range("O6").Value = 0.06415146660411865
$ws.Range("P6").Value = 0.06415146660411865
$ws.Range("Q6").Value = 19.42443556855245
$ws.Range("R6").Value = 174.819920116972
$ws.Range("S6").Value = 0.02959487777308916
$ws.Range("T6").Value = 0.02959487777308916
$ws.Range("I7").Value = 0.461328155686921
$ws.Range("J7").Value = 0.4613281556869209
$ws.Range("O7").Value = 0.3979101621202897
$ws.Range("P7").Value = 0.3979101621202898
$ws.Range("S7").Value = 0.183567161220037
$ws.Range("T7").Value = 0.183567161220037
$ws.Range("I8").Value = 0.461328155686921
$ws.Range("J8").Value = 0.4613281556869209
$ws.Range("M8").Value = 21.90816333333333
$ws.Range("N8").Value = 65.72449
$ws.Range("O8").Value = 0.1593353362087987
$ws.Range("P8").Value = 0.1593353362087987
$ws.Range("Q8").Value = 48.24517872803779
$ws.Range("R8").Value = 434.2066085523401
$ws.Range("S8").Value = 0.07350587678896058
$ws.Range("T8").Value = 0.07350587678896056
$ws.Range("I9").Value = 0.461328155686921
$ws.Range("J9").Value = 0.4613281556869209
$ws.Range("M9").Value = 52.056859
$ws.Range("N9").Value = 156.170577
$ws.Range("O9").Value = 0.3786030350667928
$ws.Range("P9").Value = 0.3786030350667929
$ws.Range("Q9").Value = 114.6372896834313
$ws.Range("R9").Value = 1031.735607150882
$ws.Range("S9").Value = 0.1746602399048342
$ws.Range("T9").Value = 0.1746602399048342
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.1506176666666667
$ws.Range("H10").Value = 0.451853
$ws.Range("I10").Value = 0.03155280162368235
$ws.Range("J10").Value = 0.03155280162368235
$ws.Range("M10").Value = 8.820647333333334
$ws.Range("N10").Value = 26.461942
$ws.Range("O10").Value = 0.06415146660411865
$ws.Range("P10").Value = 0.06415146660411865
$ws.Range("Q10").Value = 1.328545319836222
$ws.Range("R10").Value = 11.956907878526
$ws.Range("S10").Value = 0.002024158499628039
$ws.Range("T10").Value = 0.002024158499628039
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.1506176666666667
$ws.Range("H11").Value = 0.451853
$ws.Range("I11").Value = 0.03155280162368235
$ws.Range("J11").Value = 0.03155280162368235
$ws.Range("O11").Value = 0.3979101621202897
$ws.Range("P11").Value = 0.3979101621202898
$ws.Range("Q11").Value = 8.240523741451668
$ws.Range("R11").Value = 74.16471367306499
$ws.Range("S11").Value = 0.01255518040942879
$ws.Range("T11").Value = 0.01255518040942879
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.1506176666666667
$ws.Range("H12").Value = 0.451853
$ws.Range("I12").Value = 0.03155280162368235
$ws.Range("J12").Value = 0.03155280162368235
$ws.Range("M12").Value = 21.90816333333333
$ws.Range("N12").Value = 65.72449
$ws.Range("O12").Value = 0.1593353362087987
$ws.Range("P12").Value = 0.1593353362087987
$ws.Range("Q12").Value = 3.299756442218889
$ws.Range("R12").Value = 29.69780797997
$ws.Range("S12").Value = 0.005027476255038956
$ws.Range("T12").Value = 0.005027476255038956
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.1506176666666667
$ws.Range("H13").Value = 0.451853
$ws.Range("I13").Value = 0.03155280162368235
$ws.Range("J13").Value = 0.03155280162368235
$ws.Range("M13").Value = 52.056859
$ws.Range("N13").Value = 156.170577
$ws.Range("O13").Value = 0.3786030350667928
$ws.Range("P13").Value = 0.3786030350667929
$ws.Range("Q13").Value = 7.840682636575667
$ws.Range("R13").Value = 70.56614372918099
$ws.Range("S13").Value = 0.01194598645958657
$ws.Range("T13").Value = 0.01194598645958657
$ws.Range("G14").Value = 1.031517666666667
$ws.Range("H14").Value = 3.094553
$ws.Range("I14").Value = 0.2160919965629775
$ws.Range("J14").Value = 0.2160919965629775
$ws.Range("M14").Value = 8.820647333333334
$ws.Range("N14").Value = 26.461942
$ws.Range("O14").Value = 0.06415146660411865
$ws.Range("P14").Value = 0.06415146660411865
$ws.Range("Q14").Value = 9.098653555769555
$ws.Range("R14").Value = 81.887882001926
$ws.Range("S14").Value = 0.01386261850092717
$ws.Range("T14").Value = 0.01386261850092717
$ws.Range("G15").Value = 1.031517666666667
$ws.Range("H15").Value = 3.094553
$ws.Range("I15").Value = 0.2160919965629775
$ws.Range("J15").Value = 0.2160919965629775
$ws.Range("O15").Value = 0.3979101621202897
$ws.Range("P15").Value = 0.3979101621202898
$ws.Range("Q15").Value = 56.43591492295166
$ws.Range("R15").Value = 507.923234306565
$ws.Range("S15").Value = 0.08598520138527148
$ws.Range("T15").Value = 0.08598520138527149
$ws.Range("G16").Value = 1.031517666666667
$ws.Range("H16").Value = 3.094553
$ws.Range("I16").Value = 0.2160919965629775
$ws.Range("J16").Value = 0.2160919965629775
$ws.Range("M16").Value = 21.90816333333333
$ws.Range("N16").Value = 65.72449
$ws.Range("O16").Value = 0.1593353362087987
$ws.Range("P16").Value = 0.1593353362087987
$ws.Range("Q16").Value = 22.59865752255222
$ws.Range("R16").Value = 203.38791770297
$ws.Range("S16").Value = 0.03443109092439259
$ws.Range("T16").Value = 0.03443109092439259
$ws.Range("G17").Value = 1.031517666666667
$ws.Range("H17").Value = 3.094553
$ws.Range("I17").Value = 0.2160919965629775
$ws.Range("J17").Value = 0.2160919965629775
$ws.Range("M17").Value = 52.056859
$ws.Range("N17").Value = 156.170577
$ws.Range("O17").Value = 0.3786030350667928
$ws.Range("P17").Value = 0.3786030350667929
$ws.Range("Q17").Value = 53.69756972967566
$ws.Range("R17").Value = 483.2781275670809
$ws.Range("S17").Value = 0.08181308575238626
$ws.Range("T17").Value = 0.08181308575238627
